$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '80.567.98'
$ws.Range("E2").Value = '  +5.76%  '
$ws.Range("D3").Value = '3.232.96'
$ws.Range("E3").Value = '  +6.81%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '''214.64'
$ws.Range("E5").Value = '  +8.63%  '
$ws.Range("D6").Value = '''642.49'
$ws.Range("E6").Value = '  +3.97%  '
$ws.Range("D7").Value = '''0.276'
$ws.Range("E7").Value = '  +33.58%  '
$ws.Range("D8").Value = '''0.998'
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").Value = '''0.609'
$ws.Range("E9").Value = '  +11.05%  '
$ws.Range("D10").Value = '3.232.69'
$ws.Range("E10").Value = '  +6.95%  '
$ws.Range("D11").Value = '''0.628'
$ws.Range("E11").Value = '  +42.79%  '
$ws.Range("D12").Value = '''0.0000277'
$ws.Range("E12").Value = '  +44.52%  '
$ws.Range("E13").Value = '  +3.67%  '
$ws.Range("D14").Value = '''5.49'
$ws.Range("E14").Value = '  +5.42%  '
$ws.Range("D15").Value = '3.820.56'
$ws.Range("E15").Value = '  +6.81%  '
$ws.Range("D16").Value = '''32.91'
$ws.Range("E16").Value = '  +14.20%  '
$ws.Range("D17").Value = '80.177.64'
$ws.Range("E17").Value = '  +5.43%  '
$ws.Range("D18").Value = '3.215.43'
$ws.Range("E18").Value = '  +6.71%  '
$ws.Range("D19").Value = '''14.81'
$ws.Range("E19").Value = '  +10.37%  '
$ws.Range("D20").Value = '''3.07'
$ws.Range("E20").Value = '  +29.06%  '
$ws.Range("D21").Value = '''9.44'
$ws.Range("E21").Value = '  +6.02%  '
$ws.Range("D22").Value = '''450.25'
$ws.Range("E22").Value = '  +18.87%  '
$ws.Range("D23").Value = '''5.45'
$ws.Range("E23").Value = '  +24.41%  '
$ws.Range("D24").Value = '''4.90'
$ws.Range("E24").Value = '  +13.47%  '
$ws.Range("D25").Value = '''78.21'
$ws.Range("E25").Value = '  +8.12%  '
$ws.Range("D26").Value = '3.367.38'
$ws.Range("E26").Value = '  +5.96%  '
$ws.Range("D27").Value = '''11.04'
$ws.Range("E27").Value = '  +12.75%  '
$ws.Range("D28").Value = '''0.0000129'
$ws.Range("E28").Value = '  +19.81%  '
$ws.Range("D29").Value = '''1.00'
$ws.Range("E29").Value = '  +0.00%  '
$ws.Range("D30").Value = '''9.37'
$ws.Range("E30").Value = '  +13.67%  '
$ws.Range("D31").Value = '''0.998'
$ws.Range("E31").Value = '  -0.20%  '
$ws.Range("D32").Value = '''571.51'
$ws.Range("E32").Value = '  +16.29%  '
$ws.Range("D33").Value = '''1.55'
$ws.Range("E33").Value = '  +11.58%  '
$ws.Range("D34").Value = '''0.159'
$ws.Range("E34").Value = '  +29.09%  '
$ws.Range("E35").Value = '  +7.76%  '
$ws.Range("D36").Value = '''23.78'
$ws.Range("E36").Value = '  +15.91%  '
$ws.Range("E37").Value = '  +21.10%  '
$ws.Range("D38").Value = '''1.00'
$ws.Range("E38").Value = '  +0.08%  '
$ws.Range("D39").Value = '''0.418'
$ws.Range("E39").Value = '  +11.29%  '
$ws.Range("D40").Value = '''5.92'
$ws.Range("E40").Value = '  +15.96%  '
$ws.Range("D41").Value = '''164.17'
$ws.Range("E41").Value = '  +1.22%  '
$ws.Range("D42").Value = '''20.33'
$ws.Range("E42").Value = '  +1.41%  '
$ws.Range("D43").Value = '''194.30'
$ws.Range("E43").Value = '  +2.21%  '
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("D45").Value = '''1.86'
$ws.Range("E45").Value = '  +13.76%  '
$ws.Range("D46").Value = '''2.79'
$ws.Range("E46").Value = '  +15.63%  '
$ws.Range("E47").Value = '  +9.86%  '
$ws.Range("D48").Value = '''0.807'
$ws.Range("E48").Value = '  +3.72%  '
$ws.Range("B49").Value = 'OKB'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D49").Value = '''43.55'
$ws.Range("E49").Value = '  +6.19%  '
$ws.Range("B50").Value = 'Filecoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D50").Value = '''4.36'
$ws.Range("E50").Value = '  +13.22%  '
$ws.Range("D51").Value = '''0.653'
$ws.Range("E51").Value = '  +10.49%  '
